$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (C) to "Y" for all testcase rows so that
# all profile testcases run.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("C10").Value = "Y"
$ws.Range("C11").Value = "Y"
$ws.Range("C12").Value = "Y"
$ws.Range("C13").Value = "Y"

# Update the sheet's active selection
$ws.Range("B17").Select() | Out-Null
